# Refresh cached Market Board profit figures across all Leve tables.
# Values below come from an external price-data refresh; no formulas are
# involved (source workbook stores plain cached numbers in H:N).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    ,@(28, "H", 126216.125)
    ,@(28, "I", 167746.75)
    ,@(28, "J", 1624.25)
    ,@(28, "K", 167746.75)
    ,@(28, "L", 1624.25)
    ,@(28, "M", -167261.75)
    ,@(28, "N", -2594.25)
    ,@(58, "H", 2818.4)
    ,@(58, "I", 2818.4)
    ,@(58, "K", 8455.200000000001)
    ,@(58, "M", -8305.200000000001)
    ,@(80, "H", 50006470)
    ,@(80, "I", 200000500)
    ,@(80, "J", 8464.532999999999)
    ,@(80, "K", 600001500)
    ,@(80, "L", 25393.599)
    ,@(80, "M", -600000502)
    ,@(80, "N", -27389.599)
    ,@(83, "H", 50006470)
    ,@(83, "I", 200000500)
    ,@(83, "J", 8464.532999999999)
    ,@(83, "K", 1800004500)
    ,@(83, "L", 76180.79699999999)
    ,@(83, "M", -1799999508)
    ,@(83, "N", -86164.79699999999)
    ,@(88, "H", 4048.25)
    ,@(88, "I", 3846.5)
    ,@(88, "K", 3846.5)
    ,@(88, "M", -3440.5)
    ,@(91, "H", 4048.25)
    ,@(91, "I", 3846.5)
    ,@(91, "K", 3846.5)
    ,@(91, "M", -2442.5)
    ,@(98, "H", 1270)
    ,@(98, "I", 1203.2106)
    ,@(98, "K", 1203.2106)
    ,@(98, "M", 294.7893999999999)
    ,@(122, "H", 1270)
    ,@(122, "I", 1203.2106)
    ,@(122, "K", 3609.6318)
    ,@(122, "M", -1159.6318)
    ,@(132, "H", 3268.0544)
    ,@(132, "I", 3291.537)
    ,@(132, "K", 9874.610999999999)
    ,@(132, "M", -7344.610999999999)
    ,@(135, "H", 1083.2222)
    ,@(135, "I", 1138.32)
    ,@(135, "J", 394.5)
    ,@(135, "K", 10244.88)
    ,@(135, "L", 3550.5)
    ,@(135, "M", -7709.879999999999)
    ,@(135, "N", -8620.5)
    ,@(137, "H", 43677.8)
    ,@(137, "I", 51005.94)
    ,@(137, "J", 2151.6667)
    ,@(137, "K", 153017.82)
    ,@(137, "L", 6455.000100000001)
    ,@(137, "M", -150467.82)
    ,@(137, "N", -11555.0001)
    ,@(138, "H", 3159.3232)
    ,@(138, "I", 1142.2916)
    ,@(138, "J", 3804.7734)
    ,@(138, "K", 3426.8748)
    ,@(138, "L", 11414.3202)
    ,@(138, "M", 1713.1252)
    ,@(138, "N", -21694.3202)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    ,@(32, "H", 11468637)
    ,@(32, "I", 11699088)
    ,@(32, "J", 8933671)
    ,@(32, "K", 11699088)
    ,@(32, "L", 8933671)
    ,@(32, "M", -11698801)
    ,@(32, "N", -8934245)
    ,@(38, "H", 10000)
    ,@(38, "I", 0)
    ,@(38, "J", 10000)
    ,@(38, "K", 0)
    ,@(38, "L", 10000)
    ,@(38, "M", $null)
    ,@(38, "N", -10934)
    ,@(45, "H", 3482.6538)
    ,@(45, "I", 3049.2104)
    ,@(45, "J", 4659.143)
    ,@(45, "K", 3049.2104)
    ,@(45, "L", 4659.143)
    ,@(45, "M", -2672.2104)
    ,@(45, "N", -5413.143)
    ,@(61, "H", 3345.238)
    ,@(61, "I", 3473.516)
    ,@(61, "J", 2983.7273)
    ,@(61, "K", 3473.516)
    ,@(61, "L", 2983.7273)
    ,@(61, "M", -3261.516)
    ,@(61, "N", -3407.7273)
    ,@(74, "H", 2698.675)
    ,@(74, "I", 2218.6875)
    ,@(74, "K", 2218.6875)
    ,@(74, "M", -1344.6875)
    ,@(77, "H", 2698.675)
    ,@(77, "I", 2218.6875)
    ,@(77, "K", 11093.4375)
    ,@(77, "M", -6725.4375)
    ,@(102, "H", 2358.5)
    ,@(102, "I", 1082)
    ,@(102, "J", 5337)
    ,@(102, "K", 1082)
    ,@(102, "L", 5337)
    ,@(102, "M", 540)
    ,@(102, "N", -8581)
    ,@(124, "H", 30713)
    ,@(124, "J", 30713)
    ,@(124, "L", 30713)
    ,@(124, "N", -40533)
    ,@(132, "H", 2768.8865)
    ,@(132, "I", 2293.4849)
    ,@(132, "J", 4195.091)
    ,@(132, "K", 6880.4547)
    ,@(132, "L", 12585.273)
    ,@(132, "M", -4350.4547)
    ,@(132, "N", -17645.273)
    ,@(136, "H", 3345.238)
    ,@(136, "I", 3473.516)
    ,@(136, "J", 2983.7273)
    ,@(136, "K", 10420.548)
    ,@(136, "L", 8951.1819)
    ,@(136, "M", -7870.548000000001)
    ,@(136, "N", -14051.1819)
    ,@(139, "H", 81347.39999999999)
    ,@(139, "J", 81934.625)
    ,@(139, "L", 81934.625)
    ,@(139, "N", -92214.625)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    ,@(20, "H", 1589.8214)
    ,@(20, "I", 1209.6)
    ,@(20, "K", 1209.6)
    ,@(20, "M", -962.5999999999999)
    ,@(81, "H", 94267.71000000001)
    ,@(81, "J", 94267.71000000001)
    ,@(81, "L", 94267.71000000001)
    ,@(81, "N", -96389.71000000001)
    ,@(84, "H", 94267.71000000001)
    ,@(84, "J", 94267.71000000001)
    ,@(84, "L", 282803.13)
    ,@(84, "N", -293411.13)
    ,@(86, "H", 2375.111)
    ,@(86, "I", 1899)
    ,@(86, "K", 1899)
    ,@(86, "M", -776)
    ,@(89, "H", 2375.111)
    ,@(89, "I", 1899)
    ,@(89, "K", 9495)
    ,@(89, "M", -3879)
    ,@(134, "H", 2167107)
    ,@(134, "I", 3969820)
    ,@(134, "K", 11909460)
    ,@(134, "M", -11906925)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    ,@(132, "H", 6153.857)
    ,@(132, "I", 6173.385)
    ,@(132, "J", 5900)
    ,@(132, "K", 18520.155)
    ,@(132, "L", 17700)
    ,@(132, "M", -15990.155)
    ,@(132, "N", -22760)
    ,@(134, "H", 2134.3)
    ,@(134, "I", 2063.818)
    ,@(134, "J", 2328.125)
    ,@(134, "K", 6191.454000000001)
    ,@(134, "L", 6984.375)
    ,@(134, "M", -3656.454000000001)
    ,@(134, "N", -12054.375)
    ,@(140, "H", 403333)
    ,@(140, "J", 549999.5)
    ,@(140, "L", 549999.5)
    ,@(140, "N", -560359.5)
    ,@(141, "H", 419278.84)
    ,@(141, "J", 419278.84)
    ,@(141, "L", 419278.84)
    ,@(141, "N", -429638.84)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    ,@(107, "H", 934.9)
    ,@(107, "I", 1422.75)
    ,@(107, "J", 757.5)
    ,@(107, "K", 4268.25)
    ,@(107, "L", 2272.5)
    ,@(107, "M", -2348.25)
    ,@(107, "N", -6112.5)
    ,@(122, "H", 996799.1)
    ,@(122, "J", 2324065.2)
    ,@(122, "L", 20916586.8)
    ,@(122, "N", -20921486.8)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    ,@(3, "H", 4877250)
    ,@(3, "I", 5800800)
    ,@(3, "K", 5800800)
    ,@(3, "M", -5800684)
    ,@(101, "H", 69696)
    ,@(101, "J", 69696)
    ,@(101, "L", 69696)
    ,@(101, "N", -76186)
    ,@(107, "H", 1535.2222)
    ,@(107, "I", 1063.75)
    ,@(107, "J", 1912.4)
    ,@(107, "K", 1063.75)
    ,@(107, "L", 1912.4)
    ,@(107, "M", 856.25)
    ,@(107, "N", -5752.4)
    ,@(122, "H", 4162.25)
    ,@(122, "I", 4580)
    ,@(122, "J", 2909)
    ,@(122, "K", 13740)
    ,@(122, "L", 8727)
    ,@(122, "M", -11290)
    ,@(122, "N", -13627)
    ,@(126, "H", 3062.375)
    ,@(126, "I", 4300)
    ,@(126, "J", 2649.8333)
    ,@(126, "K", 12900)
    ,@(126, "L", 7949.499899999999)
    ,@(126, "M", -10430)
    ,@(126, "N", -12889.4999)
    ,@(132, "H", 3937.0588)
    ,@(132, "I", 3795.5)
    ,@(132, "K", 11386.5)
    ,@(132, "M", -8856.5)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    ,@(7, "H", 7742.619)
    ,@(7, "I", 8836)
    ,@(7, "K", 8836)
    ,@(7, "M", -8724)
    ,@(38, "H", 60055)
    ,@(38, "J", 60055)
    ,@(38, "L", 60055)
    ,@(38, "N", -60875)
    ,@(93, "H", 38462800)
    ,@(93, "I", 62500684)
    ,@(93, "J", 2182.9)
    ,@(93, "K", 62500684)
    ,@(93, "L", 2182.9)
    ,@(93, "M", -62499436)
    ,@(93, "N", -4678.9)
    ,@(126, "H", 7742.619)
    ,@(126, "I", 8836)
    ,@(126, "K", 26508)
    ,@(126, "M", -24038)
    ,@(132, "H", 45346.105)
    ,@(132, "I", 55420.684)
    ,@(132, "J", 8406)
    ,@(132, "K", 166262.052)
    ,@(132, "L", 25218)
    ,@(132, "M", -163732.052)
    ,@(132, "N", -30278)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    ,@(103, "H", 98000)
    ,@(103, "J", 98000)
    ,@(103, "L", 98000)
    ,@(103, "N", -100344)
    ,@(126, "H", 4781.875)
    ,@(126, "I", 4895.5386)
    ,@(126, "J", 4289.3335)
    ,@(126, "K", 14686.6158)
    ,@(126, "L", 12868.0005)
    ,@(126, "M", -12216.6158)
    ,@(126, "N", -17808.0005)
    ,@(132, "H", 3061.04)
    ,@(132, "I", 2422.6924)
    ,@(132, "K", 7268.0772)
    ,@(132, "M", -4738.0772)
)
foreach ($u in $updates) {
    $r = $u[0]; $colLetter = $u[1]; $val = $u[2]
    $cell = $ws.Range("$colLetter$r")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}
